$d = $word.ActiveDocument

# --- 1. Rebuild the 3rd paragraph ("-Dépassement de stack : ??") -----------
# The original paragraph holds a hidden "_GoBack" bookmark created by Word's
# last-edit tracking. Deleting the whole paragraph removes that bookmark;
# we then re-insert the paragraph (plus the two new paragraphs that follow
# it) using raw WordprocessingML so we can reproduce the exact proofErr /
# gramStart / gramEnd markup from the target revision.

$p3 = $d.Paragraphs.Item(3)
$p3.Range.Delete()

# NB: collapsing the preceding paragraph's own Range to its end lands right
# before that paragraph's pilcrow, and InsertXML there overwrites the
# paragraph instead of adding a new one after it. Using the end of the
# document's Content (the position right after the last paragraph mark)
# inserts the new paragraphs correctly.
$insertionPoint = $d.Range($d.Content.End, $d.Content.End)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newXml = @"
<w:p $w>
  <w:r>
    <w:t xml:space="preserve">-D&#233;passement de </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>stack</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t> </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>: ??</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
<w:p $w/>
<w:p $w>
  <w:r>
    <w:t>Quand un composant a besoin d&#8217;un composant du layer d&#8217;en dessous, il utilise le &#171; </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Layer.h</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t> &#187; associ&#233;.</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> Les composants d&#8217;un </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>level</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> ne doivent pas inclure le &#171; </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Layer.h</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t> &#187; de leur level.</w:t>
  </w:r>
</w:p>
"@

$insertionPoint.InsertXML($newXml)
